# feat: add 2022-Q3 data
#
# Inserts a new "2022-Q3" tab (fund-holdings detail) right after "总计" and
# before the existing "2022-Q2" tab, and adds a matching "2022-Q3" row at the
# top of the "总计" (summary) sheet's data table, shifting the older rows
# down by one.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) New "2022-Q3" worksheet: duplicate the "2022-Q2" sheet (same headers /
#    styles / column layout) and place the copy immediately before it, then
#    rename and overwrite its data cells with the 2022-Q3 figures.
# ---------------------------------------------------------------------
$q2 = $wb.Worksheets.Item("2022-Q2")
$q2.Copy($q2)
$q3 = $wb.Worksheets.Item(2)
$q3.Name = "2022-Q3"

# Row 2 - 000906 广发全球精选股票（QDII）美元现汇
$q3.Range("B2").Value = "'000906"
$q3.Range("C2").Value = "广发全球精选股票（QDII）美元现汇"
$q3.Range("D2").Value = "'21.88"
$q3.Range("E2").Value = "'79.27"
$q3.Range("F2").Value = "'4.64"
$q3.Range("G2").Value = "'1.0152"
$q3.Range("H2").Value = 6

# Row 3 - 270023 广发全球精选股票（QDII）
$q3.Range("B3").Value = "'270023"
$q3.Range("C3").Value = "广发全球精选股票（QDII）"
$q3.Range("D3").Value = "'21.88"
$q3.Range("E3").Value = "'79.27"
$q3.Range("F3").Value = "'4.64"
$q3.Range("G3").Value = "'1.0152"
$q3.Range("H3").Value = 6

# Row 4 - 001481 华宝油气（QDII）美元
$q3.Range("B4").Value = "'001481"
$q3.Range("C4").Value = "华宝油气（QDII）美元"
$q3.Range("D4").Value = "'45.98"
$q3.Range("E4").Value = "'94.53"
$q3.Range("F4").Value = "'2.05"
$q3.Range("G4").Value = "'0.9426"
$q3.Range("H4").Value = 9

# Row 5 - 162411 华宝油气（QDII）人民币A
$q3.Range("B5").Value = "'162411"
$q3.Range("C5").Value = "华宝油气（QDII）人民币A"
$q3.Range("D5").Value = "'28.25"
$q3.Range("E5").Value = "'94.53"
$q3.Range("F5").Value = "'2.05"
$q3.Range("G5").Value = "'0.5791"
$q3.Range("H5").Value = 9

# Row 6 - 007844 华宝油气（QDII）人民币 C
$q3.Range("B6").Value = "'007844"
$q3.Range("C6").Value = "华宝油气（QDII）人民币 C"
$q3.Range("D6").Value = "'17.73"
$q3.Range("E6").Value = "'94.53"
$q3.Range("F6").Value = "'2.05"
$q3.Range("G6").Value = "'0.3635"
$q3.Range("H6").Value = 9

# ---------------------------------------------------------------------
# 2) "总计" summary sheet: push existing data rows (2..6) down to (3..7)
#    and write the new 2022-Q3 row into the now-vacant row 2.
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

for ($r = 6; $r -ge 2; $r--) {
    $dest = $r + 1
    $total.Range("A$dest").Value = $r - 1
    $total.Range("B$dest").Value = $total.Range("B$r").Value2
    $total.Range("C$dest").Value = $total.Range("C$r").Value2
    $total.Range("D$dest").Value = $total.Range("D$r").Value2
}

# Row 7 (2020-Q4) is brand-new territory for column A - copy its number
# formatting/style down from the row above so it matches the rest of the
# column instead of staying default-styled.
$total.Range("A6").Copy()
$total.Range("A7").PasteSpecial(-4122)
$total.Range("A7").Value = 5

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q3"
$total.Range("C2").Value = 5
$total.Range("D2").Value = 3.92
